$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new row above the old row 4 ("Check defaults" / data rows shift
# down by one). This is the new "Button Next" control row for the new
# radiobutton control described in the commit message.
# ---------------------------------------------------------------------------
$ws.Rows("4:4").Insert()
$ws.Rows("4:4").ClearFormats()
$ws.Range("B4:I4").Clear()

$ws.Range("A4").Value = "Button Next"
$ws.Range("J4").Value = "X"

# Give the new "X" marker cell (J4) the same plain text format the other
# header/id cells use (no fill) by copying it from A1.
$ws.Range("A1").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Row 1 (the id= header row): drop the (redundant, no-fill) alternate style
# so every header cell shares the same plain style as A1.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("B1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Row 2 (VIA_Pixel9Pro_API35 input row) and row 5 (old row 4, the
# "Check defaults" input row): keep their yellow highlight fill, it just
# collapses onto the renumbered style slot once the unused duplicate style is
# gone - copy the fill forward from itself is a no-op here, nothing to do
# value/format-wise, format already yellow highlighted.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Row 6 (old row 5, the sample data row): switch its cells from the
# yellow-highlight style to the plain style used elsewhere in the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("B6:J6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Selection / cursor moves back up to J1.
# ---------------------------------------------------------------------------
$ws.Range("J1").Select()
